# Update "Pais" (countries) dashboard sheet with refreshed COVID stats and
# re-sort two pairs of rows whose totals changed order (Uganda/Liberia and
# Timor Oriental/Santa Lucia), plus bump the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp in the title row ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Agosto de 2020 a las 12:19"

# --- Row 35 (Oman) ---
$ws.Range("B35").Value = 80713
$ws.Range("C35").Value = 427
$ws.Range("D35").Value = 70910
$ws.Range("E35").Value = 9311
$ws.Range("G35").Value = 4
$ws.Range("H35").Value = 492

# --- Row 42 (Bielorrusia) ---
$ws.Range("B42").Value = 68503
$ws.Range("C42").Value = 127
$ws.Range("D42").Value = 63756
$ws.Range("E42").Value = 4167
$ws.Range("G42").Value = 3
$ws.Range("H42").Value = 580

# --- Row 44 (Rumania) ---
$ws.Range("B44").Value = 57895
$ws.Range("C44").Value = 1345
$ws.Range("D44").Value = 28992
$ws.Range("E44").Value = 26337
$ws.Range("G44").Value = 45
$ws.Range("H44").Value = 2566

# --- Row 55 (Ghana) ---
$ws.Range("B55").Value = 39642
$ws.Range("C55").Value = 567
$ws.Range("D55").Value = 36384
$ws.Range("E55").Value = 3059

# --- Row 58 (Suiza) ---
$ws.Range("B58").Value = 36108
$ws.Range("C58").Value = 181
$ws.Range("E58").Value = 2523
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 1985

# --- Row 86 (Noruega) ---
$ws.Range("D86").Value = 8857
$ws.Range("E86").Value = 296

# --- Row 88 (Malasia) ---
$ws.Range("B88").Value = 9038
$ws.Range("C88").Value = 15
$ws.Range("D88").Value = 8713
$ws.Range("E88").Value = 200

# --- Row 126 (Eslovenia) ---
$ws.Range("B126").Value = 2223
$ws.Range("C126").Value = 15
$ws.Range("D126").Value = 1909
$ws.Range("E126").Value = 190

# --- Rows 143/144: Uganda overtakes Liberia in total cases, so they swap ---
# places in the (descending-sorted) table. Row 143 becomes Uganda's
# (updated) data; row 144 becomes Liberia's (unchanged) data.
$ws.Range("A143").Value = "Uganda"
$ws.Range("B143").Value = 1223
$ws.Range("C143").Value = 10
$ws.Range("D143").Value = 1102
$ws.Range("E143").Value = 116
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 5

$ws.Range("A144").Value = "Liberia"
$ws.Range("B144").Value = 1221
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 699
$ws.Range("E144").Value = 444
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 78

# --- Row 159 (Lesoto) ---
$ws.Range("B159").Value = 742
$ws.Range("C159").Value = 16
$ws.Range("D159").Value = 175
$ws.Range("E159").Value = 544
$ws.Range("G159").Value = 2
$ws.Range("H159").Value = 23

# --- Row 160 (Vietnam) ---
$ws.Range("B160").Value = 718
$ws.Range("C160").Value = 1
$ws.Range("G160").Value = 2
$ws.Range("H160").Value = 10

# --- Rows 202/203: Timor Oriental and Santa Lucia have equal totals, but
# swap order (Timor Oriental now listed first) ---
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"

# --- Row 210 (Groenlandia) ---
$ws.Range("D210").Value = 14
$ws.Range("E210").Value = 0
